$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: Location ---
$ws.Range("B1").Value = "Location"
$ws.Range("B2").Value = "Chennai"
$ws.Range("B3").Value = "Pune"
$ws.Range("B4").Value = "Delhi"

# --- Column C: Car Name ---
$ws.Range("C1").Value = "Car Name"
$ws.Range("C2").Value = "Hyundai I10"
$ws.Range("C3").Value = "Maruti Swift"
$ws.Range("C4").Value = "Hyundai Santro Xing"

# --- Column D: Email ---
$ws.Range("D2").Value = "hhasj"

# --- Hyperlink on D3 (email) ---
$ws.Range("D3").Value = "S@n.com"
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:S@n.com")

$ws.Range("D4").Value = "ma.com"
$ws.Range("D1").Value = "Email"

# --- Column widths (closest achievable to source 18.90625 / 17.81640625
#     after the host's internal character->pixel snapping) ---
$ws.Columns.Item(3).ColumnWidth = 18.0
$ws.Columns.Item(4).ColumnWidth = 17.0

# --- Selection ---
$ws.Range("D1").Select() | Out-Null
